$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Ano 2021)
$ws.Range("B3").Value = 956164.54
$ws.Range("C3").Value = 726.5777332171136
$ws.Range("D3").Value = 2159
$ws.Range("E3").Value = 2159
$ws.Range("F3").Value = 442.8738026864289
$ws.Range("G3").Value = -1.224152306616322

# Row 4 (Ano 2022)
$ws.Range("B4").Value = 1772367.68
$ws.Range("C4").Value = 85.36220554675663
$ws.Range("D4").Value = 2707
$ws.Range("E4").Value = 2707
$ws.Range("F4").Value = 654.7350129294422
$ws.Range("G4").Value = 47.83782850958536

# Row 5 (Ano 2023)
$ws.Range("B5").Value = 2845202.68
$ws.Range("C5").Value = 60.53117601422298
$ws.Range("D5").Value = 3334
$ws.Range("E5").Value = 3334
$ws.Range("F5").Value = 853.3901259748051
$ws.Range("G5").Value = 30.34129978119424

# Row 6 (Ano 2024)
$ws.Range("B6").Value = 4447357.64
$ws.Range("C6").Value = 56.31074971432262
$ws.Range("D6").Value = 4736
$ws.Range("E6").Value = 4736
$ws.Range("F6").Value = 939.0535557432431
$ws.Range("G6").Value = 10.03801510716884

# Row 7 (Ano 2025)
$ws.Range("B7").Value = 2272603.46
$ws.Range("C7").Value = -48.89991667051989
$ws.Range("D7").Value = 2283
$ws.Range("E7").Value = 2283
$ws.Range("F7").Value = 995.4461060008761
$ws.Range("G7").Value = 6.005253897686313
